$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of this weekly block (row 589), shifting
# all existing data (rows 589-621) down by 3 rows.
$ws.Rows("589:591").Insert()

# Populate the 3 newly inserted rows with this week's data
# (Fruta, Terminal La Palmera de La Serena, Mango - Coquimbo), one row per
# calidad (Especial / Primera / Segunda).
$calidades = @("Especial", "Primera", "Segunda")

for ($i = 0; $i -lt 3; $i++) {
    $r = 589 + $i
    $ws.Cells.Item($r, 1).Value = 8
    $ws.Cells.Item($r, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = 44610
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100108
    $ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($r, 9).Value = 100108002
    $ws.Cells.Item($r, 10).Value = "Mango"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = $calidades[$i]
    $ws.Cells.Item($r, 13).Value = 512
    $ws.Cells.Item($r, 14).Value = 6000
    $ws.Cells.Item($r, 15).Value = 6500
    $ws.Cells.Item($r, 16).Value = 6250
    $ws.Cells.Item($r, 17).Value = "`$/bandeja 4 kilos"
    $ws.Cells.Item($r, 18).Value = "Perú"
    $ws.Cells.Item($r, 19).Value = 1562
    $ws.Cells.Item($r, 20).Value = 4
}
